# Update "Inscritos" counts on the "Inscricoes" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 22
$ws.Range("E12").Value = 21
$ws.Range("E14").Value = 31
$ws.Range("E16").Value = 263
$ws.Range("E18").Value = 74
